# Update numeric values in Sheet1 to reflect revised algorithm output
# (commit message: "Update Name of Algo")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D8").Value  = -7.575
$ws.Range("D10").Value = -7.999
$ws.Range("D12").Value = -7.937
$ws.Range("E13").Value = 12.583
$ws.Range("D18").Value = -8.038999999999998
$ws.Range("D25").Value = -8.183
